$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Invalid block: Unexpected tag EOF missing [ENDFOR]",
    $false,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Invalid block: Unexpected tag EOF missing [ENDFOR] while parsing m:for v | self.eClassifiers",
    2
)
